# Updates cryptos list cell values to match the latest scraped data.
# D column holds price text (often formatted like "1.234.56"); E column holds
# padded percentage-change text (e.g. "  +1.61%  "). Both are stored as plain
# text in the workbook, so numeric-looking D values are written into cells that
# are first switched to the "@" (Text) number format to stop Excel from silently
# re-interpreting them as numbers and dropping meaningful digits (e.g. "1.00").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) values that are unambiguous text already (contain more
# than one "." so Excel will not coerce them to a number) ---
$ws.Range("D2").Value = "36.236.79"
$ws.Range("D3").Value = "2.020.82"
$ws.Range("D14").Value = "2.313.24"
$ws.Range("D16").Value = "2.021.99"
$ws.Range("D18").Value = "36.399.94"
$ws.Range("D46").Value = "1.299.60"
$ws.Range("D49").Value = "2.204.44"

# --- Price (column D) values that look like plain numbers; force text format
# first so the literal string (e.g. trailing zeros) is preserved ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.62"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.652"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.66"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.81"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0707"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0979"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.06"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.791"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.81"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.36"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "233.63"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.43"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.87"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.61"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.91"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.19"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.26"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0570"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0872"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.92"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.10"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.835"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.30"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "94.69"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.78"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.30"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0812"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.14"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.73"

# --- Volume/percentage-change (column E) values ---
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("E3").Value = "  +6.69%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  -5.94%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.28%  "
$ws.Range("E9").Value = "  +5.00%  "
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("E11").Value = "  -5.99%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("E16").Value = "  +6.44%  "
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("E19").Value = "  -4.47%  "
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("E21").Value = "  -5.20%  "
$ws.Range("E22").Value = "  -4.53%  "
$ws.Range("E23").Value = "  -7.00%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -9.54%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("E28").Value = "  +6.75%  "
$ws.Range("E29").Value = "  -10.87%  "
$ws.Range("E30").Value = "  -6.32%  "
$ws.Range("E31").Value = "  +48.82%  "
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("E33").Value = "  -6.23%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("E36").Value = "  +18.67%  "
$ws.Range("E37").Value = "  -8.01%  "
$ws.Range("E38").Value = "  +7.44%  "
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("E40").Value = "  -12.29%  "
$ws.Range("E41").Value = "  -4.55%  "
$ws.Range("E42").Value = "  -8.06%  "
$ws.Range("E43").Value = "  +15.87%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  -10.50%  "
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  +6.36%  "
$ws.Range("E50").Value = "  -8.92%  "
$ws.Range("E51").Value = "  +12.75%  "

